# Update cryptocurrency price (D) and 1h volume-change (E) columns
# to reflect the refreshed figures from the Feb 16 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.908.88'
$ws.Range("E2").Value = '  +0.43%  '

$ws.Range("D3").Value = '2.787.00'
$ws.Range("E3").Value = '  -1.43%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '358.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.51%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.565'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.71%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.596'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0854'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.02%  '

$ws.Range("E12").Value = '  +0.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.31%  '

$ws.Range("D15").Value = '3.218.54'
$ws.Range("E15").Value = '  -1.82%  '

$ws.Range("D16").Value = '2.808.12'
$ws.Range("E16").Value = '  -0.91%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.935'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.27%  '

$ws.Range("D18").Value = '51.822.01'
$ws.Range("E18").Value = '  +0.45%  '

$ws.Range("E19").Value = '  +1.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.96%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.83%  '

$ws.Range("D22").Value = '0.0₃0980'
$ws.Range("E22").Value = '  -1.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '274.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.86%  '

$ws.Range("E25").Value = '  -0.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("E28").Value = '  +1.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.145'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.72%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0466'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.44%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '51.52'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.71%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.70'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.84%  '

$ws.Range("E35").Value = '  +2.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.25'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.01%  '

$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.24'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.68%  '

$ws.Range("E39").Value = '  +0.33%  '

$ws.Range("E40").Value = '  -3.23%  '

$ws.Range("E41").Value = '  +2.36%  '

$ws.Range("E42").Value = '  -1.81%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '122.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.24%  '

$ws.Range("E44").Value = '  -2.36%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.02%  '

$ws.Range("D46").Value = '2.073.80'
$ws.Range("E46").Value = '  -0.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.53%  '

$ws.Range("E49").Value = '  +0.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.935'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.83%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.96'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.65%  '
